$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.924.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.655.88"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.41%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.45"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "189.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.619"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.63%  "
$ws.Range("E8").Value = "  -0.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.699"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.15%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "57.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.84%  "
$ws.Range("E11").Value = "  -6.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000273"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -6.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.16"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.240.71"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.661.22"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.74%  "
$ws.Range("E16").Value = "  +0.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.84"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.13%  "
$ws.Range("E18").Value = "  -1.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "67.737.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "399.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "87.55"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.93%  "
$ws.Range("E25").Value = "  -3.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.83%  "
$ws.Range("E27").Value = "  -0.31%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.79"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.24"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "31.79"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.34"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "44.80"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "66.31"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.116"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.14%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "608.52"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.23%  "
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("E38").Value = "  -1.41%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0773"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -13.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.134"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.89"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0425"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.54"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.135"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.794.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.36%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.15"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "142.85"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.04%  "
$ws.Range("E50").Value = "  -3.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.49"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -16.28%  "
